$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Sales Return"

# Header row
$ws.Range("A1").Value = "Sl"
$ws.Range("B1").Value = "Module Name"

# Data rows (entered in the order the shared-string table was built:
# confirm was typed in last, after edit/delete, even though it sits
# visually above them in the sheet)
$ws.Range("B2").Value = "sales_return"
$ws.Range("B3").Value = "sales_return_list"
$ws.Range("B4").Value = "sales_return_view"
$ws.Range("B6").Value = "sales_return_edit"
$ws.Range("B7").Value = "sales_return_delete"
$ws.Range("B5").Value = "sales_return_confirm"

# Column widths (closest the engine's column-width quantization allows
# to the authored 12.42578125 / 31.140625 "characters" widths)
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 30.333333333333336

# Selection on next empty cell below data
$ws.Range("B8").Select()
